# Updated main GSC export data:
# The GSC export was re-pulled and the very first day in the "Chart" sheet
# (2025-10-08 -- a day with no recorded video-indexing data) dropped out of
# the reporting window. Remove that obsolete first data row; Excel shifts
# every following row up by one, renumbers the shared-string table, and
# updates the sheet's used-range dimension automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 1 is the header (Date / No video indexed / Video indexed / Impressions).
# Row 2 is the obsolete 2025-10-08 entry being dropped from the export.
$ws.Rows.Item(2).Delete()
